$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-20 Friday", 2)

$d.Content.Find.Execute("51÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=", 2)
$d.Content.Find.Execute("11÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷8=", 2)
$d.Content.Find.Execute("87÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷4=", 2)
$d.Content.Find.Execute("63÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=", 2)
$d.Content.Find.Execute("91÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=", 2)
$d.Content.Find.Execute("66÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=", 2)
$d.Content.Find.Execute("79÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=", 2)
$d.Content.Find.Execute("15÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=", 2)
$d.Content.Find.Execute("10÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷9=", 2)
$d.Content.Find.Execute("94÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷4=", 2)
$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=", 2)
$d.Content.Find.Execute("52÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷8=", 2)
$d.Content.Find.Execute("35÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷6=", 2)
$d.Content.Find.Execute("82÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷6=", 2)
$d.Content.Find.Execute("32÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=", 2)
$d.Content.Find.Execute("58÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=", 2)
$d.Content.Find.Execute("75÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=", 2)
$d.Content.Find.Execute("82÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=", 2)
$d.Content.Find.Execute("20÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=", 2)
$d.Content.Find.Execute("52÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=", 2)
$d.Content.Find.Execute("71÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=", 2)
$d.Content.Find.Execute("77÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=", 2)
$d.Content.Find.Execute("35÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷8=", 2)
$d.Content.Find.Execute("39÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=", 2)
$d.Content.Find.Execute("94÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷3=", 2)
